# Generate Report for Handoff
#
# A new source file (cdcd26a4-5d27-49de-b574-042e8a82d80f.md) was handed off.
# Its status is recorded as a brand-new row inserted above the existing
# "ef9f47ca-34b1-4ebe-bbfe-af8189f04721" entry on every sheet (Overview,
# zh-cn, de-de), pushing the older entry down one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Make room for the new entry above the existing data row.
$ws1.Rows.Item(2).Insert()

# The hyperlink that used to live on A2 logically belongs to the file that
# just moved down to row 3, so drop the stale link and recreate it there.
$ws1.Range("A2").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/946a9afd840cf7b605bac3083b13e214fe3a734a/e2e/ef9f47ca-34b1-4ebe-bbfe-af8189f04721.md", "", "", "ef9f47ca-34b1-4ebe-bbfe-af8189f04721.md")

# Populate the new row with the newly handed-off file.
$ws1.Range("A2").Value = "cdcd26a4-5d27-49de-b574-042e8a82d80f.md"
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"
$ws1.Range("D2").Value = "2016-30-12 16:30:17"
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/946a9afd840cf7b605bac3083b13e214fe3a734a/e2e/cdcd26a4-5d27-49de-b574-042e8a82d80f.md", "", "", "cdcd26a4-5d27-49de-b574-042e8a82d80f.md")

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Rows.Item(2).Insert()

$ws2.Range("A2").Hyperlinks.Delete()
$ws2.Range("B2").Hyperlinks.Delete()
$ws2.Range("D2").Hyperlinks.Delete()

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/946a9afd840cf7b605bac3083b13e214fe3a734a/e2e/ef9f47ca-34b1-4ebe-bbfe-af8189f04721.md", "", "", "ef9f47ca-34b1-4ebe-bbfe-af8189f04721.md")
$ws2.Hyperlinks.Add($ws2.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/946a9afd840cf7b605bac3083b13e214fe3a734a/e2e/ef9f47ca-34b1-4ebe-bbfe-af8189f04721.md", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f27def8e3288f9514f4bbd2e7af586fe6d199e4e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ef9f47ca-34b1-4ebe-bbfe-af8189f04721.922f61ca9dbc7456b5825222b5f5d0450a2f55c6.zh-cn.xlf", "", "", "ef9f47ca-34b1-4ebe-bbfe-af8189f04721.922f61ca9dbc7456b5825222b5f5d0450a2f55c6.zh-cn.xlf")

$ws2.Range("A2").Value = "cdcd26a4-5d27-49de-b574-042e8a82d80f.md"
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("D2").Value = "cdcd26a4-5d27-49de-b574-042e8a82d80f.47f99b9c6100d8f91e24bc41387943de4db3c6d9.zh-cn.xlf"
$ws2.Range("E2").Value = "2016-03-12 16:30:14"
$ws2.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H2").Value = "0001-01-01 00:00:00"
$ws2.Range("I2").Value = "Include"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/946a9afd840cf7b605bac3083b13e214fe3a734a/e2e/cdcd26a4-5d27-49de-b574-042e8a82d80f.md", "", "", "cdcd26a4-5d27-49de-b574-042e8a82d80f.md")
$ws2.Hyperlinks.Add($ws2.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/946a9afd840cf7b605bac3083b13e214fe3a734a/e2e/cdcd26a4-5d27-49de-b574-042e8a82d80f.md", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f27def8e3288f9514f4bbd2e7af586fe6d199e4e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/cdcd26a4-5d27-49de-b574-042e8a82d80f.47f99b9c6100d8f91e24bc41387943de4db3c6d9.zh-cn.xlf", "", "", "cdcd26a4-5d27-49de-b574-042e8a82d80f.47f99b9c6100d8f91e24bc41387943de4db3c6d9.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Rows.Item(2).Insert()

$ws3.Range("A2").Hyperlinks.Delete()
$ws3.Range("B2").Hyperlinks.Delete()
$ws3.Range("D2").Hyperlinks.Delete()

$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/946a9afd840cf7b605bac3083b13e214fe3a734a/e2e/ef9f47ca-34b1-4ebe-bbfe-af8189f04721.md", "", "", "ef9f47ca-34b1-4ebe-bbfe-af8189f04721.md")
$ws3.Hyperlinks.Add($ws3.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/946a9afd840cf7b605bac3083b13e214fe3a734a/e2e/ef9f47ca-34b1-4ebe-bbfe-af8189f04721.md", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/43dbbf109e92f9e29c69f14f94057dfee208b74d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ef9f47ca-34b1-4ebe-bbfe-af8189f04721.922f61ca9dbc7456b5825222b5f5d0450a2f55c6.de-de.xlf", "", "", "ef9f47ca-34b1-4ebe-bbfe-af8189f04721.922f61ca9dbc7456b5825222b5f5d0450a2f55c6.de-de.xlf")

$ws3.Range("A2").Value = "cdcd26a4-5d27-49de-b574-042e8a82d80f.md"
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("D2").Value = "cdcd26a4-5d27-49de-b574-042e8a82d80f.47f99b9c6100d8f91e24bc41387943de4db3c6d9.de-de.xlf"
$ws3.Range("E2").Value = "2016-03-12 16:30:17"
$ws3.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H2").Value = "0001-01-01 00:00:00"
$ws3.Range("I2").Value = "Include"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/946a9afd840cf7b605bac3083b13e214fe3a734a/e2e/cdcd26a4-5d27-49de-b574-042e8a82d80f.md", "", "", "cdcd26a4-5d27-49de-b574-042e8a82d80f.md")
$ws3.Hyperlinks.Add($ws3.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/946a9afd840cf7b605bac3083b13e214fe3a734a/e2e/cdcd26a4-5d27-49de-b574-042e8a82d80f.md", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/43dbbf109e92f9e29c69f14f94057dfee208b74d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/cdcd26a4-5d27-49de-b574-042e8a82d80f.47f99b9c6100d8f91e24bc41387943de4db3c6d9.de-de.xlf", "", "", "cdcd26a4-5d27-49de-b574-042e8a82d80f.47f99b9c6100d8f91e24bc41387943de4db3c6d9.de-de.xlf")

Write-Host "Generated report for handoff: added cdcd26a4-5d27-49de-b574-042e8a82d80f.md row"
